# Add 2022-Q3 data sheet and row in the summary ("总计") sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q3" worksheet, positioned right after "总计" ---
$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

$newSheet = $wb.Worksheets.Add($q1Sheet)
$newSheet.Name = "2022-Q3"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").Style = "Accent1"

# Data rows
$data = @(
    @(0, "006165", "建信中证1000指数增强A", "3.87", "84.02", "1.75", "0.0677", 1),
    @(1, "006166", "建信中证1000指数增强C", "1.89", "84.02", "1.75", "0.0331", 1),
    @(2, "013442", "建信中证1000指数增强E", "0.18", "84.02", "1.75", "0.0032", 1)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- 2. Insert a new row for "2022-Q3" in the "总计" (summary) sheet ---
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.1

# Re-number the index column (A) for the remaining rows since they shifted
for ($i = 3; $i -le 8; $i++) {
    $totalSheet.Cells.Item($i, 1).Value = $i - 2
}
